# Applies the "LinuxForHealth" re-branding edit described in the commit:
#   - URL, Version, Date and Publisher metadata values updated on the
#     "Metadata" sheet
#   - the stray ele-1/ext-1 constraint duplicated onto the top-level
#     "Extension" row (AI2) on the "Elements" sheet is cleared (it still
#     correctly appears on the "Extension.extension" row, AI4)

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/recorder-facility"
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
$elements.Range("AI2").ClearContents()
